# This script updates the "Starting_State" transition-probability matrix on
# Sheet1 of the workbook. The sheet holds an 18x18 matrix of game-state
# transition probabilities (rows/cols 0..17, labelled Af0..Br0). Several
# cells that previously held a placeholder 0 now hold the actual simulated
# transition probabilities (as decimal fractions that sum to 1 across each
# row), reflecting results from running more games through the simulator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (state "Af0")
$ws.Range("B2").Value = 0.09090909090909091
$ws.Range("C2").Value = 0.6363636363636364
$ws.Range("P2").Value = 0.2727272727272727

# Row 3 (state "Af1")
$ws.Range("B3").Value = 0.2857142857142857
$ws.Range("P3").Value = 0.5714285714285714
$ws.Range("S3").Value = 0.1428571428571428

# Row 4 (state "Af2")
$ws.Range("S4").Value = 1

# Row 6 (state "Ai0")
$ws.Range("B6").Value = 0.2
$ws.Range("D6").Value = 0.2
$ws.Range("O6").Value = 0.2
$ws.Range("Q6").Value = 0.2
$ws.Range("S6").Value = 0.2

# Row 7 (state "Ai1")
$ws.Range("B7").Value = 0.2
$ws.Range("J7").Value = 0.4
$ws.Range("R7").Value = 0.2
$ws.Range("S7").Value = 0.2

# Row 8 (state "Ai2")
$ws.Range("B8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.06666666666666667
$ws.Range("Q8").Value = 0.3333333333333333
$ws.Range("R8").Value = 0.06666666666666667
$ws.Range("S8").Value = 0.4666666666666667

# Row 9 (state "Ai3")
$ws.Range("F9").Value = 0.5
$ws.Range("J9").Value = 0.5

# Row 10 (state "Ar0")
$ws.Range("B10").Value = 0.06944444444444445
$ws.Range("E10").Value = 0.01388888888888889
$ws.Range("F10").Value = 0.02777777777777778
$ws.Range("J10").Value = 0.1388888888888889
$ws.Range("O10").Value = 0.01388888888888889
$ws.Range("Q10").Value = 0.3333333333333333
$ws.Range("R10").Value = 0.1388888888888889
$ws.Range("S10").Value = 0.2638888888888889

# Row 11 (state "Bf0")
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("J11").Value = 0.1428571428571428
$ws.Range("K11").Value = 0.1428571428571428
$ws.Range("L11").Value = 0.5714285714285714

# Row 12 (state "Bf1")
$ws.Range("G12").Value = 1

# Row 15 (state "Bi0")
$ws.Range("H15").Value = 0.1111111111111111
$ws.Range("J15").Value = 0.4444444444444444
$ws.Range("K15").Value = 0.2222222222222222
$ws.Range("S15").Value = 0.2222222222222222

# Row 16 (state "Bi1")
$ws.Range("H16").Value = 0.2857142857142857
$ws.Range("J16").Value = 0.7142857142857143

# Row 17 (state "Bi2")
$ws.Range("H17").Value = 0.09677419354838709
$ws.Range("J17").Value = 0.6451612903225806
$ws.Range("K17").Value = 0.06451612903225806
$ws.Range("O17").Value = 0.03225806451612903
$ws.Range("S17").Value = 0.1612903225806452

# Row 18 (state "Bi3")
$ws.Range("H18").Value = 0.1818181818181818
$ws.Range("J18").Value = 0.7272727272727273
$ws.Range("K18").Value = 0.09090909090909091

# Row 19 (state "Br0")
$ws.Range("F19").Value = 0.02380952380952381
$ws.Range("H19").Value = 0.1666666666666667
$ws.Range("I19").Value = 0.04761904761904762
$ws.Range("J19").Value = 0.5
$ws.Range("K19").Value = 0.02380952380952381
$ws.Range("O19").Value = 0.119047619047619
$ws.Range("S19").Value = 0.119047619047619
